$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "pair_kind" value (generic) for the 4 practice rows ---
$ws.Range("J2").Value = "generic"
$ws.Range("J3").Value = "generic"
$ws.Range("J4").Value = "generic"
$ws.Range("J5").Value = "generic"

# --- New "stim details" block added at the bottom of the sheet ---
$ws.Range("A27").Value = "stim details"

# Header row for the new block
$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

# Data rows: month number (A) + word_type (B)
$stimRows = @(
    @{ Row = 29; Month = 6; WordType = "video" },
    @{ Row = 30; Month = 6; WordType = "video" },
    @{ Row = 31; Month = 7; WordType = "video" },
    @{ Row = 32; Month = 7; WordType = "video" },
    @{ Row = 33; Month = 6; WordType = "audio" },
    @{ Row = 34; Month = 6; WordType = "audio" },
    @{ Row = 35; Month = 7; WordType = "audio" },
    @{ Row = 36; Month = 7; WordType = "audio" }
)

foreach ($r in $stimRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Month
    $ws.Cells.Item($r.Row, 2).Value = $r.WordType
}
